$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front (shifts B:F existing data right)
$ws.Columns.Item(1).Insert()

# Remove old hyperlinks (they still reference the pre-shift column D and
# would otherwise remain anchored there) and re-add them against the new
# column E locations, preserving original order/targets.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:a@bc.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:a@bc.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:f@fl.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:d@ef.stu")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:a@bc.com")

# Populate the new "Username" column
$ws.Range("A1").Value = "Username"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2").Value = "aaa"
$ws.Range("A3").Value = "bbb"
$ws.Range("A4").Value = "ccc"
$ws.Range("A5").Value = "ddd"
$ws.Range("A6").Value = "eee"

# Move the active selection to A7, matching the author's final cursor spot
$ws.Range("A7").Select()
